$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("X2:AB2").ClearContents()
$ws.Range("AD2").ClearContents()
